# Certification Data for Data Analytics
# Inserts a new "Data Analytics Specialty" row into the certification sheet
# directly above the existing "Database Specialty" row (new row 7), pushing
# the remaining rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting rows 7:14 down to 8:15.
$ws.Rows.Item(7).Insert(-4121, 0)

# The plain Insert() leaves the new row's cells with a border-less default
# style, so copy formatting down from the row above (row 6) to match the
# other data rows (s="5"/"6", with the bottom border + wrap text on col D).
$ws.Range("A6:J6").Copy()
$ws.Range("A7:J7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new row's data.
$ws.Range("A7").Value = "ff44f4b1-a3de-40a2-bbfc-2f1d53b60c86"
$ws.Range("B7").Value = "Data Analytics Specialty"
$ws.Range("C7").Value = "L622V9QJEJ41135K"
$ws.Range("D7").Value = "[`n""Ability to define AWS data analytics services and understand how they integrate with each other"",`n""Ability to explain how AWS data analytics services fit in the data lifecycle of collection, storage, processing, and visualization""`n]"
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 2024
$ws.Range("G7").Value = 8
$ws.Range("H7").Value = 2021
$ws.Range("I7").Value = "SPECIALTY"
$ws.Range("J7").Value = "AWS"

# Keep the row height consistent with the other data rows (15pt) instead of
# letting the wrapped description text auto-expand it.
$ws.Rows.Item(7).RowHeight = 15

# Update the selection to match the post-edit state (activeCell E16).
$ws.Range("E16").Select()
